$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

foreach ($r in 2..8) {
    $ws.Cells.Item($r, 3).Value = 45185
}
